# Added code for daily challenge day 13
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DSA Tracker")
$ws2 = $wb.Worksheets.Item("Practice Previous")

# Colors (OLE BGR-packed ints), matching the fills already used in the sheet
$colorEasy   = 5287936   # 00B050 - "Easy"
$colorMedium = 65535     # FFFF00 - "Medium"
$colorHard   = 255       # FF0000 - "Hard"

# Use an existing, already-date-formatted cell as the formatting template so
# the pasted-in cells land on the very same style record instead of minting
# a fresh (but equivalent) one for every single cell.
$dateTemplate = $ws2.Range("I2")

function Set-Difficulty($ws, $rowNum, $level, $dateSerial) {
    $hCell = $ws.Cells.Item($rowNum, 8)
    $iCell = $ws.Cells.Item($rowNum, 9)

    $hCell.Value = $level
    if ($level -eq "Easy") {
        $hCell.Interior.Color = $colorEasy
    } elseif ($level -eq "Medium") {
        $hCell.Interior.Color = $colorMedium
    } else {
        $hCell.Interior.Color = $colorHard
    }

    $dateTemplate.Copy()
    $iCell.PasteSpecial(-4122) | Out-Null
    $iCell.Value = $dateSerial
}

# ---------------------------------------------------------------
# Sheet "DSA Tracker": fill in Difficulty / Next Revision for rows 2-13
# ---------------------------------------------------------------
Set-Difficulty $ws1 2  "Easy"   46065
Set-Difficulty $ws1 3  "Easy"   46065
Set-Difficulty $ws1 4  "Easy"   46065
Set-Difficulty $ws1 5  "Easy"   46065
Set-Difficulty $ws1 6  "Easy"   46065
Set-Difficulty $ws1 7  "Easy"   46065
Set-Difficulty $ws1 8  "Medium" 46040
Set-Difficulty $ws1 9  "Hard"   46036
Set-Difficulty $ws1 10 "Easy"   46065
Set-Difficulty $ws1 11 "Easy"   46065
Set-Difficulty $ws1 12 "Easy"   46065
Set-Difficulty $ws1 13 "Easy"   46065

# Row 46's revision date moves two days out
$ws1.Cells.Item(46, 9).Value = 46037

# ---------------------------------------------------------------
# Sheet "Practice Previous": refresh difficulty / revision date on rows 3-10
# ---------------------------------------------------------------
Set-Difficulty $ws2 3  "Medium" 46041
Set-Difficulty $ws2 4  "Medium" 46041
Set-Difficulty $ws2 5  "Medium" 46041
Set-Difficulty $ws2 6  "Medium" 46041
Set-Difficulty $ws2 7  "Hard"   46037
Set-Difficulty $ws2 8  "Medium" 46041
Set-Difficulty $ws2 9  "Medium" 46041
Set-Difficulty $ws2 10 "Medium" 46041

# D5 / D7 become real hyperlinks (reusing the look of the other linked cells)
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://leetcode.com/problems/sort-colors/submissions/1881574577/") | Out-Null
$ws2.Range("D2").Copy()
$ws2.Range("D5").PasteSpecial(-4122) | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("D7"), "https://leetcode.com/problems/maximum-subarray/") | Out-Null
$ws2.Range("D2").Copy()
$ws2.Range("D7").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------
# New rows: Daily Challenge day 13 additions
# ---------------------------------------------------------------
# Column E (Data Structure) carries a shaded-fill style on every existing row;
# reuse that exact style (instead of re-deriving the theme colour) via paste.
$colETemplate = $ws2.Range("E2")

$ws2.Cells.Item(11, 1).Value = 10
$ws2.Cells.Item(11, 2).Value = "Next Permutation"
$ws2.Cells.Item(11, 3).Value = 31
$ws2.Cells.Item(11, 4).Value = "https://leetcode.com/problems/next-permutation/"
$colETemplate.Copy()
$ws2.Cells.Item(11, 5).PasteSpecial(-4122) | Out-Null
$ws2.Cells.Item(11, 5).Value = "Array"
$ws2.Cells.Item(11, 6).Value = "Traverse, Reverse"
$ws2.Cells.Item(11, 7).Value = "Leetcode"
Set-Difficulty $ws2 11 "Hard" 46037

$ws2.Cells.Item(12, 1).Value = 11
$ws2.Cells.Item(12, 2).Value = "Longest Consecutive Sequence"
$ws2.Cells.Item(12, 3).Value = 128
$ws2.Cells.Item(12, 4).Value = "https://leetcode.com/problems/longest-consecutive-sequence/"
$colETemplate.Copy()
$ws2.Cells.Item(12, 5).PasteSpecial(-4122) | Out-Null
$ws2.Cells.Item(12, 5).Value = "Array"
$ws2.Cells.Item(12, 6).Value = "Traverse, Hashing"
$ws2.Cells.Item(12, 7).Value = "Leetcode"
Set-Difficulty $ws2 12 "Medium" 46040

$ws2.Cells.Item(13, 1).Value = 12
$ws2.Cells.Item(13, 2).Value = "Pascal's triangle"
$ws2.Cells.Item(13, 3).Value = 118
$ws2.Cells.Item(13, 4).Value = "https://leetcode.com/problems/pascals-triangle/description/"
$colETemplate.Copy()
$ws2.Cells.Item(13, 5).PasteSpecial(-4122) | Out-Null
$ws2.Cells.Item(13, 5).Value = "Array"
$ws2.Cells.Item(13, 6).Value = "Traverse"
$ws2.Cells.Item(13, 7).Value = "Leetcode"
Set-Difficulty $ws2 13 "Medium" 46040

# Cosmetic: selection the author ended up on after editing
$ws1.Activate()
$ws1.Range("M53").Select() | Out-Null
$ws2.Range("H23").Select() | Out-Null

Write-Host "edit complete"
